$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 657.6
$ws.Range("I2").Value = 737.5
$ws.Range("J2").Value = 604.3333
$ws.Range("K2").Value = 737.5
$ws.Range("L2").Value = 604.3333
$ws.Range("M2").Value = -624.5
$ws.Range("N2").Value = -830.3333
$ws.Range("H9").Value = 235
$ws.Range("I9").Value = 82.5
$ws.Range("J9").Value = 540
$ws.Range("K9").Value = 82.5
$ws.Range("L9").Value = 540
$ws.Range("M9").Value = 86.5
$ws.Range("N9").Value = -878
$ws.Range("I40").Value = 1199.1666
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 1199.1666
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -1024.1666
$ws.Range("N40").Value = -1550
$ws.Range("H92").Value = 670.6429000000001
$ws.Range("I92").Value = 766.1111
$ws.Range("K92").Value = 766.1111
$ws.Range("M92").Value = 481.8889
$ws.Range("H116").Value = 2985.5
$ws.Range("I116").Value = 2985.5
$ws.Range("K116").Value = 2985.5
$ws.Range("M116").Value = 456.5
$ws.Range("H129").Value = 1100
$ws.Range("I129").Value = 1100
$ws.Range("K129").Value = 3300
$ws.Range("M129").Value = 1700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6294.7617
$ws.Range("I32").Value = 6294.7617
$ws.Range("K32").Value = 6294.7617
$ws.Range("M32").Value = -6007.7617
$ws.Range("H35").Value = 1200
$ws.Range("I35").Value = 1200
$ws.Range("K35").Value = 1200
$ws.Range("M35").Value = -794
$ws.Range("H45").Value = 2326
$ws.Range("I45").Value = 1042.4
$ws.Range("J45").Value = 3242.8572
$ws.Range("K45").Value = 1042.4
$ws.Range("L45").Value = 3242.8572
$ws.Range("M45").Value = -665.4000000000001
$ws.Range("N45").Value = -3996.8572
$ws.Range("H70").Value = 99999
$ws.Range("J70").Value = 99999
$ws.Range("L70").Value = 99999
$ws.Range("N70").Value = -100539
$ws.Range("H73").Value = 99999
$ws.Range("J73").Value = 99999
$ws.Range("L73").Value = 99999
$ws.Range("N73").Value = -101871
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 2000
$ws.Range("M74").Value = -1126
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 10000
$ws.Range("M77").Value = -5632
$ws.Range("H97").Value = 1364.2
$ws.Range("I97").Value = 1124.2222
$ws.Range("J97").Value = 1724.1666
$ws.Range("K97").Value = 1124.2222
$ws.Range("L97").Value = 1724.1666
$ws.Range("M97").Value = -628.2221999999999
$ws.Range("N97").Value = -2716.1666
$ws.Range("H102").Value = 5903.3335
$ws.Range("I102").Value = 5903.3335
$ws.Range("K102").Value = 5903.3335
$ws.Range("M102").Value = -4281.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10999.667
$ws.Range("J86").Value = 12999
$ws.Range("L86").Value = 12999
$ws.Range("N86").Value = -15245
$ws.Range("H89").Value = 10999.667
$ws.Range("J89").Value = 12999
$ws.Range("L89").Value = 64995
$ws.Range("N89").Value = -76227
$ws.Range("H106").Value = 5221.5
$ws.Range("J106").Value = 5221.5
$ws.Range("L106").Value = 5221.5
$ws.Range("N106").Value = -7745.5
$ws.Range("H134").Value = 345
$ws.Range("I134").Value = 345
$ws.Range("K134").Value = 1035
$ws.Range("M134").Value = 1500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 313.1111
$ws.Range("I23").Value = 335.4
$ws.Range("J23").Value = 285.25
$ws.Range("K23").Value = 1006.2
$ws.Range("L23").Value = 855.75
$ws.Range("M23").Value = -771.1999999999999
$ws.Range("N23").Value = -1325.75
$ws.Range("H97").Value = 487.375
$ws.Range("I97").Value = 514
$ws.Range("J97").Value = 443
$ws.Range("K97").Value = 1542
$ws.Range("L97").Value = 1329
$ws.Range("M97").Value = -1046
$ws.Range("N97").Value = -2321
$ws.Range("H117").Value = 2175.2727
$ws.Range("J117").Value = 2671.5
$ws.Range("L117").Value = 8014.5
$ws.Range("N117").Value = -14898.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 50000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H16").Value = 50000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H62").Value = 46000
$ws.Range("I62").Value = 46000
$ws.Range("K62").Value = 46000
$ws.Range("M62").Value = -45314
$ws.Range("H65").Value = 46000
$ws.Range("I65").Value = 46000
$ws.Range("K65").Value = 138000
$ws.Range("M65").Value = -134568
$ws.Range("H132").Value = 4115.3076
$ws.Range("I132").Value = 3874.9167
$ws.Range("K132").Value = 11624.7501
$ws.Range("M132").Value = -9094.750100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1975
$ws.Range("I61").Value = 1975
$ws.Range("K61").Value = 1975
$ws.Range("M61").Value = -1773
$ws.Range("H63").Value = 46000
$ws.Range("I63").Value = 46000
$ws.Range("K63").Value = 46000
$ws.Range("M63").Value = -45251
$ws.Range("H66").Value = 46000
$ws.Range("I66").Value = 46000
$ws.Range("K66").Value = 138000
$ws.Range("M66").Value = -134256
$ws.Range("H113").Value = 1975
$ws.Range("I113").Value = 1975
$ws.Range("K113").Value = 1975
$ws.Range("M113").Value = 195

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 959.6667
$ws.Range("I107").Value = 439.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1318.5
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 601.5
$ws.Range("N107").Value = -9840
$ws.Range("H122").Value = 2161.5715
$ws.Range("I122").Value = 2122.6667
$ws.Range("J122").Value = 2395
$ws.Range("K122").Value = 6368.000100000001
$ws.Range("L122").Value = 7185
$ws.Range("M122").Value = -3918.000100000001
$ws.Range("N122").Value = -12085
